# Update mods data [2025-12-29 15:11:27]
# Append a new data row (row 49) to the ModCounts sheet:
#   A49 = "2025/12/29"  (text, same style as the existing date column)
#   B49 = "逃离鸭科夫"   (text)
#   C49 = 1111           (number)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 48
$newRow  = $lastRow + 1

# Clone the formatting (style) of the previous row onto the new row first,
# so the new cells pick up the same centered style used throughout the
# table (style index "1" in the sheet) without us having to rebuild it by
# hand.
$ws.Range("A" + $lastRow + ":C" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":C" + $newRow).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Write the date as a literal formula-text first (so it is never
# auto-parsed into a date serial number), then flatten it down to a plain
# static value while keeping the formatting already applied above.
$ws.Range("A" + $newRow).Formula = "=""2025/12/29"""
$ws.Range("A" + $newRow).Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

$ws.Range("B" + $newRow).Value = "逃离鸭科夫"
$ws.Range("C" + $newRow).Value = 1111
